# Updates the cryptocurrency price/volume snapshot in Sheet1 (columns B-E)
# per the "Updated cryptos list" GitHub Actions commit.
#
# Price cells (column D) that read as plain decimal numbers (e.g. "1.00",
# "0.390") are written through Range.Formula with a leading apostrophe so
# Excel keeps them as literal text (matching the source data, which stores
# prices like "0.999" / "1.00" / "3.387.75" as text, not numbers) instead of
# auto-converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.676.51'
$ws.Range("E2").Value = '  +1.19%  '

$ws.Range("D3").Value = '3.393.61'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Formula = "'577.62"
$ws.Range("E5").Value = '  +1.03%  '

$ws.Range("D6").Formula = "'137.60"
$ws.Range("E6").Value = '  +1.52%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '3.392.38'
$ws.Range("E8").Value = '  +0.81%  '

$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").Formula = "'7.51"
$ws.Range("E10").Value = '  -1.08%  '

$ws.Range("E11").Value = '  +2.65%  '

$ws.Range("D12").Formula = "'0.390"
$ws.Range("E12").Value = '  +0.42%  '

$ws.Range("D13").Value = '3.970.09'
$ws.Range("E13").Value = '  +0.58%  '

$ws.Range("E14").Value = '  +1.85%  '

$ws.Range("E15").Value = '  +2.39%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Formula = "'25.95"
$ws.Range("E16").Value = '  +2.80%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.385.95'
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("D18").Value = '61.736.31'
$ws.Range("E18").Value = '  +0.90%  '

$ws.Range("D19").Formula = "'14.25"
$ws.Range("E19").Value = '  +2.33%  '

$ws.Range("E20").Value = '  +0.64%  '

$ws.Range("E21").Value = '  +0.10%  '

$ws.Range("D22").Formula = "'377.71"
$ws.Range("E22").Value = '  +1.33%  '

$ws.Range("D23").Formula = "'0.559"
$ws.Range("E23").Value = '  -1.51%  '

$ws.Range("D24").Value = '3.524.92'
$ws.Range("E24").Value = '  +0.61%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").Value = '  +8.09%  '

$ws.Range("D27").Formula = "'71.18"
$ws.Range("E27").Value = '  +0.98%  '

$ws.Range("D28").Formula = "'1.67"
$ws.Range("E28").Value = '  -0.06%  '

$ws.Range("E29").Value = '  -2.03%  '

$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("E31").Value = '  +3.96%  '

$ws.Range("D32").Formula = "'8.23"
$ws.Range("E32").Value = '  +1.12%  '

$ws.Range("D33").Formula = "'2.18"
$ws.Range("E33").Value = '  +1.34%  '

$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("D35").Formula = "'23.43"
$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").Formula = "'5.34"
$ws.Range("E36").Value = '  -4.06%  '

$ws.Range("E37").Value = '  +0.48%  '

$ws.Range("D38").Formula = "'6.85"
$ws.Range("E38").Value = '  -1.32%  '

$ws.Range("D39").Formula = "'164.89"
$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("D40").Formula = "'0.0784"
$ws.Range("E40").Value = '  -0.22%  '

$ws.Range("E41").Value = '  +2.62%  '

$ws.Range("D42").Formula = "'0.782"
$ws.Range("E42").Value = '  +2.84%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Formula = "'1.74"
$ws.Range("E43").Value = '  +8.78%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Formula = "'0.999"
$ws.Range("E44").Value = '  -0.19%  '

$ws.Range("D45").Formula = "'25.32"
$ws.Range("E45").Value = '  +9.33%  '

$ws.Range("D46").Formula = "'4.42"
$ws.Range("E46").Value = '  +0.36%  '

$ws.Range("D47").Formula = "'41.46"
$ws.Range("E47").Value = '  +0.35%  '

$ws.Range("D48").Formula = "'6.87"
$ws.Range("E48").Value = '  -1.11%  '

$ws.Range("D49").Formula = "'22.77"
$ws.Range("E49").Value = '  -1.84%  '

$ws.Range("D50").Value = '2.350.32'
$ws.Range("E50").Value = '  +6.04%  '

$ws.Range("E51").Value = '  +1.94%  '
